# Updates cryptos list values (price + 1h volume change) per row,
# mirroring a scheduled data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.299.63"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "2.654.22"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'605.87"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "'153.06"
$ws.Range("E6").Value = "  +5.59%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("D10").Value = "'0.387"
$ws.Range("E10").Value = "  +6.31%  "

$ws.Range("D11").Value = "'5.62"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").Value = "'28.20"
$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").Value = "3.130.39"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "64.131.19"
$ws.Range("E15").Value = "  +1.50%  "

$ws.Range("D16").Value = "'0.0000148"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").Value = "2.648.22"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  +6.63%  "

$ws.Range("E19").Value = "  +4.18%  "

$ws.Range("D20").Value = "'347.96"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "'6.92"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'5.56"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'66.68"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("E25").Value = "  +13.32%  "

$ws.Range("E26").Value = "  +8.05%  "

$ws.Range("E27").Value = "  +3.60%  "

$ws.Range("D28").Value = "'8.18"
$ws.Range("E28").Value = "  +3.18%  "

$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").Value = "'549.08"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("E33").Value = "  +6.23%  "

$ws.Range("D34").Value = "'1.78"
$ws.Range("E34").Value = "  -1.54%  "

$ws.Range("D35").Value = "'5.35"
$ws.Range("E35").Value = "  +4.33%  "

$ws.Range("D36").Value = "'168.33"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +6.46%  "

$ws.Range("D40").Value = "'19.41"
$ws.Range("E40").Value = "  +1.46%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "'166.32"
$ws.Range("E42").Value = "  -3.49%  "

$ws.Range("D43").Value = "'40.08"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").Value = "'3.86"
$ws.Range("E44").Value = "  +2.37%  "

$ws.Range("D45").Value = "'0.0577"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "'21.91"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("D47").Value = "'0.631"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("E48").Value = "  +14.93%  "

$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").Value = "'0.0967"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").Value = "'19.07"
$ws.Range("E51").Value = "  +1.42%  "
